$wb = $excel.ActiveWorkbook

# Give the three generic tabs more descriptive names.
$wb.Worksheets.Item("Sheet1").Name = "namedTab1"
$wb.Worksheets.Item("Sheet2").Name = "namedTab2"
$wb.Worksheets.Item("Sheet3").Name = "namedTab3"

# The work being reviewed lives on the third tab - make it the active/selected
# sheet (this flips which sheetView carries tabSelected="1" and updates the
# workbook's bookViews/activeTab).
$wb.Worksheets.Item("namedTab3").Activate()
